$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 16.94798469543457
$ws.Range("C3").Value = 16.83807373046875
$ws.Range("C4").Value = 23.12397956848145
$ws.Range("C5").Value = 17.36974716186523
$ws.Range("C6").Value = 16.25609397888184
